# Apply the LOM3236 disciplina content refresh:
# - fills in the previously-missing Objetivos/Programa resumido/Programa/Bibliografia
#   bodies and the second docente row, restoring the A/B/C row alignment for rows 13-25.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 10
$ws.Range("B10").Value = 'Desenvolver conhecimento teórico e prático dos processos de fabricação de equipamentos e dispositivos requeridos para o desenvolvimento de produtos e protótipos. Conhecer os requisitos e efeitos dos processos de fabricação de forma a permitir, interagir, criar e executar projetos ao longo de sua vida profissional.'
$ws.Range("C10").Value = 'Desenvolver conhecimento teórico e prático dos processos de fabricação de equipamentos e dispositivos requeridos para o desenvolvimento de produtos e protótipos. Conhecer os requisitos e efeitos dos processos de fabricação de forma a permitir, interagir, criar e executar projetos ao longo de sua vida profissional.'

# Row 13
$ws.Range("A13").ClearContents()
$ws.Range("B13").Value = '519033 - Carlos Yujiro Shigue'
$ws.Range("C13").Value = '519033 - Carlos Yujiro Shigue'
$ws.Rows.Item(13).EntireRow.AutoFit()

# Row 14
$ws.Range("A14").ClearContents()
$ws.Range("B14").Value = '5817692 - Katia Cristiane Gandolpho Candioto'
$ws.Range("C14").Value = '5817692 - Katia Cristiane Gandolpho Candioto'
$ws.Rows.Item(14).EntireRow.AutoFit()

# Row 15
$ws.Range("A15").Value = 'Programa resumido:'
$ws.Range("B15").Value = 'Introdução aos processos de fabricação. Processos de união de materiais. Revisão de projeto assistido por computador (CAD). Manufatura auxiliada por computador (CAM). Sistemas de produção flexíveis. Prototipação rápida.'
$ws.Range("C15").Value = 'Introdução aos processos de fabricação. Processos de união de materiais. Revisão de projeto assistido por computador (CAD). Manufatura auxiliada por computador (CAM). Sistemas de produção flexíveis. Prototipação rápida.'
$ws.Rows.Item(15).RowHeight = 60

# Row 16
$ws.Range("A16").Value = 'Short syllabus:'
$ws.Range("B16").Value = 'Introduction to manufacturing processes. Material joining processes. Computer-aided design (CAD) review. Computer Aided Manufacturing (CAM). Flexible production systems. Rapid prototyping.'
$ws.Range("C16").Value = 'Introduction to manufacturing processes. Material joining processes. Computer-aided design (CAD) review. Computer Aided Manufacturing (CAM). Flexible production systems. Rapid prototyping.'
$ws.Rows.Item(16).RowHeight = 60

# Row 17
$ws.Range("A17").Value = 'Programa:'
$ws.Range("B17").Value = 'Classificação dos processos de fabricação. Fundição. Metalurgia do pó. Usinagem: processos, fundamentos e condições econômicas. Máquinas-ferramentas. Conformação mecânica. Processos de união de materiais. Manufatura auxiliada por computador (CAM). Linguagens de programação para controle numérico. Máquinas ferramentas de controle numérico. Sequência de fabricação de produtos. Noções de automação dos processos de manufatura. Prototipação rápida. Sistemas de prototipação rápida (sólido, líquido e pó).'
$ws.Range("C17").Value = 'Classificação dos processos de fabricação. Fundição. Metalurgia do pó. Usinagem: processos, fundamentos e condições econômicas. Máquinas-ferramentas. Conformação mecânica. Processos de união de materiais. Manufatura auxiliada por computador (CAM). Linguagens de programação para controle numérico. Máquinas ferramentas de controle numérico. Sequência de fabricação de produtos. Noções de automação dos processos de manufatura. Prototipação rápida. Sistemas de prototipação rápida (sólido, líquido e pó).'
$ws.Rows.Item(17).RowHeight = 120

# Row 18
$ws.Range("A18").Value = 'Syllabus:'
$ws.Range("B18").Value = 'Classification of manufacturing processes. Foundry. Powder metallurgy. Machining: processes, fundamentals and economic conditions. Machine tools. Mechanical conformation. Material joining processes. Computer Aided Manufacturing (CAM). Programming languages for numerical control. Numerical control machine tools. Product manufacturing sequence. Notions of automation of manufacturing processes. Rapid prototyping. Rapid prototyping systems (solid, liquid and powder).'
$ws.Range("C18").Value = 'Classification of manufacturing processes. Foundry. Powder metallurgy. Machining: processes, fundamentals and economic conditions. Machine tools. Mechanical conformation. Material joining processes. Computer Aided Manufacturing (CAM). Programming languages for numerical control. Numerical control machine tools. Product manufacturing sequence. Notions of automation of manufacturing processes. Rapid prototyping. Rapid prototyping systems (solid, liquid and powder).'
$ws.Rows.Item(18).RowHeight = 120

# Row 19
$ws.Range("A19").Value = 'Avaliação:'
$ws.Range("B19").ClearContents()
$ws.Range("C19").ClearContents()
$ws.Rows.Item(19).EntireRow.AutoFit()

# Row 20
$ws.Range("A20").Value = 'Método:'
$ws.Range("B20").Value = 'Aulas expositivas e aulas práticas em laboratórios de usinagem. Projeto e fabricação de protótipos.'
$ws.Range("C20").Value = 'Aulas expositivas e aulas práticas em laboratórios de usinagem. Projeto e fabricação de protótipos.'

# Row 21
$ws.Range("A21").Value = 'Critério:'
$ws.Range("B21").Value = 'Média ponderada das atividades práticas desenvolvidas, trabalhos e relatórios.'
$ws.Range("C21").Value = 'Média ponderada das atividades práticas desenvolvidas, trabalhos e relatórios.'
$ws.Rows.Item(21).RowHeight = 60

# Row 22
$ws.Range("A22").Value = 'Norma de recuperação:'
$ws.Range("B22").Value = 'Devido às características da disciplina, não será oferecida recuperação.'
$ws.Range("C22").Value = 'Devido às características da disciplina, não será oferecida recuperação.'
$ws.Rows.Item(22).RowHeight = 60

# Row 23
$ws.Range("A23").Value = 'Bibliografia:'
$ws.Range("B23").Value = @'
SMID, P. CNC Programming Handbook, Industrial Press, 2007.
GROOVER, M.; ZIMEMERS, E. Computer Aided Design and Manufacturing, Prentice-Hall, 1984.
STENERSON, J.; CURRAN, K. Computer Numerical Control: Operation and Programming, Prentice Hall, 2006.
SIMON, W. Numerical Control of Machine Tools, Edward Arnold, 1973. 
MILNER, D.; VASILOV, V.: Computer Aided Engineering for Manufacture. Kogan Page, 1986.
CHUA, C. K.; LEONG, K. F. Rapid Prototyping: Principles and Applications, World Scientific Publishing, 2010. MESSLER, R. W. Joining of Materials and Structures, Butterworth-Heinemann, 2004.
KIMINAMI, C. S.; CASTRO, W. B.; OLIVEIRA, M. F. Introdução aos processos de Fabricação de Produtos Metálicos, Blucher, 2013.
MEYERS, M.A. AND CHAWLA, K.K.; Mechanical Behavior of Materials, Prentice-Hall, Upper Saddle River-NJ (EUA), 1999.
GIESECKE, F. E. Comunicação Gráfica Moderna, Editora Bookman, 2002.
CRUZ, M. D. Catia V5r20 - Modelagem, Montagem e Detalhamento, ERICA, 2010.
FISCHER, U; GOMERINGER, R; HEINZLER, M; ET AL. Manual de Tecnologia Metal Mecânica, Blucher, 2011.
JACK, H. Projeto, Planejamento e Gestão de Produtos  Uma abordagem para engenharia, Campus-Elsevier, 2014.
SWIFT, K.G.; BOOKER, P.D. Seleção de processos de manufatura, Campus-Elsevier, 2014.
'@
$ws.Range("C23").Value = @'
SMID, P. CNC Programming Handbook, Industrial Press, 2007.
GROOVER, M.; ZIMEMERS, E. Computer Aided Design and Manufacturing, Prentice-Hall, 1984.
STENERSON, J.; CURRAN, K. Computer Numerical Control: Operation and Programming, Prentice Hall, 2006.
SIMON, W. Numerical Control of Machine Tools, Edward Arnold, 1973. 
MILNER, D.; VASILOV, V.: Computer Aided Engineering for Manufacture. Kogan Page, 1986.
CHUA, C. K.; LEONG, K. F. Rapid Prototyping: Principles and Applications, World Scientific Publishing, 2010. MESSLER, R. W. Joining of Materials and Structures, Butterworth-Heinemann, 2004.
KIMINAMI, C. S.; CASTRO, W. B.; OLIVEIRA, M. F. Introdução aos processos de Fabricação de Produtos Metálicos, Blucher, 2013.
MEYERS, M.A. AND CHAWLA, K.K.; Mechanical Behavior of Materials, Prentice-Hall, Upper Saddle River-NJ (EUA), 1999.
GIESECKE, F. E. Comunicação Gráfica Moderna, Editora Bookman, 2002.
CRUZ, M. D. Catia V5r20 - Modelagem, Montagem e Detalhamento, ERICA, 2010.
FISCHER, U; GOMERINGER, R; HEINZLER, M; ET AL. Manual de Tecnologia Metal Mecânica, Blucher, 2011.
JACK, H. Projeto, Planejamento e Gestão de Produtos  Uma abordagem para engenharia, Campus-Elsevier, 2014.
SWIFT, K.G.; BOOKER, P.D. Seleção de processos de manufatura, Campus-Elsevier, 2014.
'@
$ws.Rows.Item(23).RowHeight = 120

# Row 24
$ws.Range("A24").Value = 'Requisitos:'

# Row 25
$ws.Range("B25").Value = @'
LOM3204 -  Desenho Técnico e Projeto Assistido por Computador  (Requisito)

'@
$ws.Range("C25").Value = @'
LOM3204 -  Desenho Técnico e Projeto Assistido por Computador  (Requisito)

'@
$ws.Rows.Item(25).RowHeight = 30
